$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws_sheet1 = $wb.Worksheets.Item("展览")
$ws_sheet1.Range("F2").Value = 2319
$ws_sheet1.Range("F3").Value = 393
$ws_sheet1.Range("F4").Value = 197
$ws_sheet1.Range("F5").Value = 319
$ws_sheet1.Range("F6").Value = 319
$ws_sheet1.Range("F7").Value = 510
$ws_sheet1.Range("F9").Value = 741
$ws_sheet1.Range("F11").Value = 739
$ws_sheet1.Range("F13").Value = 80
$ws_sheet1.Range("F14").Value = 380
$ws_sheet1.Range("F15").Value = 11
$ws_sheet1.Range("F16").Value = 1002
$ws_sheet1.Range("F17").Value = 18732
$ws_sheet1.Range("G17").Value = "已售罄"
$ws_sheet1.Range("F18").Value = 459
$ws_sheet1.Range("F19").Value = 46
$ws_sheet1.Range("F20").Value = 196
$ws_sheet1.Range("F21").Value = 275
$ws_sheet1.Range("F22").Value = 164
$ws_sheet1.Range("F23").Value = 133
$ws_sheet1.Range("F26").Value = 172
$ws_sheet1.Range("F28").Value = 312
$ws_sheet1.Range("F29").Value = 127

# Sheet: 演出 (sheet2)
$ws_sheet2 = $wb.Worksheets.Item("演出")
$ws_sheet2.Range("F4").Value = 170
$ws_sheet2.Range("F5").Value = 1
$ws_sheet2.Range("F8").Value = 216
$ws_sheet2.Range("F9").Value = 3347
$ws_sheet2.Range("F11").Value = 68
$ws_sheet2.Range("F14").Value = 27
$ws_sheet2.Range("F17").Value = 2836

# Sheet: 本地生活 (sheet3)
$ws_sheet3 = $wb.Worksheets.Item("本地生活")
$ws_sheet3.Range("F2").Value = 262
$ws_sheet3.Range("F3").Value = 78
$ws_sheet3.Range("F4").Value = 533
$ws_sheet3.Range("F5").Value = 198

# Sheet: 全部类型 (sheet4)
$ws_sheet4 = $wb.Worksheets.Item("全部类型")
$ws_sheet4.Range("F2").Value = 262
$ws_sheet4.Range("F3").Value = 78
$ws_sheet4.Range("F6").Value = 2319
$ws_sheet4.Range("F7").Value = 533
$ws_sheet4.Range("F8").Value = 393
$ws_sheet4.Range("F9").Value = 197
$ws_sheet4.Range("F10").Value = 319
$ws_sheet4.Range("F11").Value = 319
$ws_sheet4.Range("F12").Value = 510
$ws_sheet4.Range("F13").Value = 170
$ws_sheet4.Range("F15").Value = 1
$ws_sheet4.Range("F18").Value = 198
$ws_sheet4.Range("F19").Value = 741
$ws_sheet4.Range("F21").Value = 739
$ws_sheet4.Range("F23").Value = 80
$ws_sheet4.Range("F24").Value = 380
$ws_sheet4.Range("F25").Value = 11
$ws_sheet4.Range("F26").Value = 1002
$ws_sheet4.Range("F27").Value = 18734
$ws_sheet4.Range("G27").Value = "暂时售罄"
$ws_sheet4.Range("F28").Value = 216
$ws_sheet4.Range("F29").Value = 3347
$ws_sheet4.Range("F31").Value = 68
$ws_sheet4.Range("F33").Value = 459
$ws_sheet4.Range("F34").Value = 46
$ws_sheet4.Range("F35").Value = 196
$ws_sheet4.Range("F37").Value = 27
$ws_sheet4.Range("F38").Value = 275
$ws_sheet4.Range("F39").Value = 164
$ws_sheet4.Range("F40").Value = 133
$ws_sheet4.Range("F45").Value = 172
$ws_sheet4.Range("F47").Value = 312
$ws_sheet4.Range("F48").Value = 127
$ws_sheet4.Range("F49").Value = 2836
